$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 1265
$ws.Range("B4").Value = "Start from intro"
$ws.Range("C4").Value = "Need to fix this with what Alicia says"

$ws.Range("C4").Select()
